# AJUSTE: archivo y enlace
#
# - Update the note in B1 describing how the price is defined
#   (old: "Precio( se definira como precio base del producto o atributo)"
#    new: "Precio se definira como el precio del rol y ciudad")
# - Add a new (empty) underlined cell at B6, mirroring the other
#   placeholder-style cells already in the sheet (D10, D16, C20)
# - Move/leave the active selection on B6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price-definition note text in B1
$ws.Range("B1").Value = "Precio se definira como el precio del rol y ciudad"

# Add the new styled (underlined) empty cell in row 6
$ws.Range("B6").Font.Underline = 2

# Move the selection to B6
$ws.Range("B6").Select()
